$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Hoja1" to "Formulario"
$ws.Name = "Formulario"

# Update row heights (in points) to match the new, more compact layout
$ws.Rows.Item(2).RowHeight = 39.75
$ws.Rows.Item(3).RowHeight = 52.5
$ws.Rows.Item(4).RowHeight = 39.75
$ws.Rows.Item(5).RowHeight = 52.5
$ws.Rows.Item(6).RowHeight = 52.5
$ws.Rows.Item(7).RowHeight = 39.75
$ws.Rows.Item(8).RowHeight = 39.75
$ws.Rows.Item(9).RowHeight = 52.5
$ws.Rows.Item(10).RowHeight = 52.5
$ws.Rows.Item(11).RowHeight = 39.75
$ws.Rows.Item(12).RowHeight = 39.75
$ws.Rows.Item(13).RowHeight = 39.75
$ws.Rows.Item(14).RowHeight = 52.5
$ws.Rows.Item(15).RowHeight = 39.75
$ws.Rows.Item(16).RowHeight = 27

# Select cell A3 so the view reflects the saved selection
$ws.Range("A3").Select()
